$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Slightly narrow the first three (label) columns.
$ws.Columns("A:C").ColumnWidth = 35.6

# Add the new "2023" year column (T), mirroring the formatting of the
# existing 2022 column (S) immediately to its left.
$ws.Range("S4").Copy($ws.Range("T4"))
$ws.Range("T4").Value = 2023

$ws.Range("S5").Copy($ws.Range("T5"))
$ws.Range("T5").Value = 40

# Reset the view back to the top-left, with a plain A1 selection instead
# of the stray scrolled/selected state left over from editing.
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("A1").Select()
